$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Style = "Normal"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.449.66'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.11%  '

# Row 3
$ws.Range("D3").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.852.22'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.15%  '

# Row 4
$ws.Range("E4").Value = '  -0.06%  '

# Row 5
$ws.Range("D5").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '242.20'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.77%  '

# Row 6
$ws.Range("D6").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6304'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.53%  '

# Row 7
$ws.Range("E7").Value = '  -0.02%  '

# Row 8
$ws.Range("D8").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07596'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.75%  '

# Row 9
$ws.Range("D9").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2983'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.52%  '

# Row 10
$ws.Range("D10").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '24.46'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.53%  '

# Row 11
$ws.Range("D11").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07708'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.93%  '

# Row 12
$ws.Range("D12").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.929.37'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.95%  '

# Row 13
$ws.Range("D13").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.009'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.94%  '

# Row 14
$ws.Range("D14").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6899'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.33%  '

# Row 15
$ws.Range("D15").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '83.49'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.24%  '

# Row 16
$ws.Range("D16").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.000009914'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.00%  '

# Row 17
$ws.Range("D17").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.183.70'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +3.55%  '

# Row 18
$ws.Range("D18").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.183'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.20%  '

# Row 19
$ws.Range("D19").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '29.603.79'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.69%  '

# Row 20
$ws.Range("D20").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '233.80'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.11%  '

# Row 21
$ws.Range("D21").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.54'
$ws.Range("D21").Style = "Normal"

# Row 22
$ws.Range("B22").Value = 'Chainlink'
$ws.Range("C22").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D22").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.711'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.78%  '

# Row 23
$ws.Range("B23").Value = 'Dai'
$ws.Range("C23").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D23").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.001'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.02%  '

# Row 25
$ws.Range("D25").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '155.44'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.84%  '

# Row 26
$ws.Range("E26").Value = '  -2.44%  '

# Row 27
$ws.Range("D27").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.477'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.92%  '

# Row 28
$ws.Range("D28").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.72'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.88%  '

# Row 29
$ws.Range("E29").Value = '  -0.87%  '

# Row 30
$ws.Range("D30").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05820'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.59%  '

# Row 31
$ws.Range("D31").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.267'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.86%  '

# Row 32
$ws.Range("D32").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.130'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.59%  '

# Row 33
$ws.Range("D33").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.021'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.71%  '

# Row 34
$ws.Range("D34").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.904'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.81%  '

# Row 35
$ws.Range("D35").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.171'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.67%  '

# Row 36
$ws.Range("D36").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7259'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.05%  '

# Row 37
$ws.Range("D37").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.590'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.81%  '

# Row 38
$ws.Range("D38").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.260.14'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +4.60%  '

# Row 39
$ws.Range("D39").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.798'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.38%  '

# Row 40
$ws.Range("D40").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01807'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.82%  '

# Row 41
$ws.Range("D41").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9098'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.36%  '

# Row 42
$ws.Range("D42").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.143'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.51%  '

# Row 43
$ws.Range("D43").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.096.27'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.72%  '

# Row 44
$ws.Range("D44").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.000'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.12%  '

# Row 45
$ws.Range("D45").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '67.86'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.10%  '

# Row 46
$ws.Range("D46").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '101.79'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.64%  '

# Row 47
$ws.Range("B47").Value = 'Aptos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D47").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '7.355'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.92%  '

# Row 48
$ws.Range("B48").Value = 'BabyDogeCoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D48").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.00000000118'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -4.00%  '

# Row 49
$ws.Range("D49").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.194'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.71%  '

# Row 50
$ws.Range("D50").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.4042'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.44%  '

# Row 51
$ws.Range("D51").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.710'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.69%  '
